# ============================================================
# 1. Append the 2025-11-12 vs MIL game row (row 12) to the four
#    per-game log sheets: Points, Assists, Rebounds, 3PM
# ============================================================
$wb = $excel.ActiveWorkbook

# --- Points ---
$ws = $wb.Worksheets.Item("Points")
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-11-12"
$ws.Range("B12").Value = "MIL"
$ws.Range("C12").Value = 20
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 16
$ws.Range("G12").Value = 16
$ws.Range("H12").Value = 17
$ws.Range("I12").Value = 7
$ws.Range("J12").Value = 11
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 13
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 2
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 5

# --- Assists ---
$ws = $wb.Worksheets.Item("Assists")
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-11-12"
$ws.Range("B12").Value = "MIL"
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 8
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 2
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 2

# --- Rebounds ---
$ws = $wb.Worksheets.Item("Rebounds")
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-11-12"
$ws.Range("B12").Value = "MIL"
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 9
$ws.Range("I12").Value = 13
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 2
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 3

# --- 3PM ---
$ws = $wb.Worksheets.Item("3PM")
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-11-12"
$ws.Range("B12").Value = "MIL"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 1

# ============================================================
# 2. Add the new "Team Points" sheet at the end of the workbook
# ============================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$teamPoints = $wb.Worksheets.Add($null, $lastSheet)
$teamPoints.Name = "Team Points"

$teamPoints.Range("A1").Value = "Game Time (PST)"
$teamPoints.Range("B1").Value = "Opponent"
$teamPoints.Range("C1").Value = "Team Points"
$teamPoints.Range("D1").Value = "Opponent Points"
$teamPoints.Range("E1").Value = "Game Total Points"

$teamPoints.Range("A2").NumberFormat = "@"
$teamPoints.Range("A2").Value = "2025-10-22"
$teamPoints.Range("B2").Value = "BKN"
$teamPoints.Range("C2").Value = 136
$teamPoints.Range("D2").Value = 117
$teamPoints.Range("E2").Value = 253
$teamPoints.Range("A3").NumberFormat = "@"
$teamPoints.Range("A3").Value = "2025-10-25"
$teamPoints.Range("B3").Value = "PHI"
$teamPoints.Range("C3").Value = 121
$teamPoints.Range("D3").Value = 125
$teamPoints.Range("E3").Value = 246
$teamPoints.Range("A4").NumberFormat = "@"
$teamPoints.Range("A4").Value = "2025-10-26"
$teamPoints.Range("B4").Value = "WAS"
$teamPoints.Range("C4").Value = 139
$teamPoints.Range("D4").Value = 113
$teamPoints.Range("E4").Value = 252
$teamPoints.Range("A5").NumberFormat = "@"
$teamPoints.Range("A5").Value = "2025-10-28"
$teamPoints.Range("B5").Value = "MIA"
$teamPoints.Range("C5").Value = 117
$teamPoints.Range("D5").Value = 144
$teamPoints.Range("E5").Value = 261
$teamPoints.Range("A6").NumberFormat = "@"
$teamPoints.Range("A6").Value = "2025-10-30"
$teamPoints.Range("B6").Value = "ORL"
$teamPoints.Range("C6").Value = 107
$teamPoints.Range("D6").Value = 123
$teamPoints.Range("E6").Value = 230
$teamPoints.Range("A7").NumberFormat = "@"
$teamPoints.Range("A7").Value = "2025-11-01"
$teamPoints.Range("B7").Value = "MIN"
$teamPoints.Range("C7").Value = 105
$teamPoints.Range("D7").Value = 122
$teamPoints.Range("E7").Value = 227
$teamPoints.Range("A8").NumberFormat = "@"
$teamPoints.Range("A8").Value = "2025-11-02"
$teamPoints.Range("B8").Value = "UTA"
$teamPoints.Range("C8").Value = 126
$teamPoints.Range("D8").Value = 103
$teamPoints.Range("E8").Value = 229
$teamPoints.Range("A9").NumberFormat = "@"
$teamPoints.Range("A9").Value = "2025-11-04"
$teamPoints.Range("B9").Value = "NOP"
$teamPoints.Range("C9").Value = 112
$teamPoints.Range("D9").Value = 116
$teamPoints.Range("E9").Value = 228
$teamPoints.Range("A10").NumberFormat = "@"
$teamPoints.Range("A10").Value = "2025-11-07"
$teamPoints.Range("B10").Value = "MIA"
$teamPoints.Range("C10").Value = 108
$teamPoints.Range("D10").Value = 126
$teamPoints.Range("E10").Value = 234
$teamPoints.Range("A11").NumberFormat = "@"
$teamPoints.Range("A11").Value = "2025-11-10"
$teamPoints.Range("B11").Value = "LAL"
$teamPoints.Range("C11").Value = 111
$teamPoints.Range("D11").Value = 121
$teamPoints.Range("E11").Value = 232
$teamPoints.Range("A12").NumberFormat = "@"
$teamPoints.Range("A12").Value = "2025-11-12"
$teamPoints.Range("B12").Value = "MIL"
$teamPoints.Range("C12").Value = 111
$teamPoints.Range("D12").Value = 100
$teamPoints.Range("E12").Value = 211

# ============================================================
# 3. Recompute + re-sort the four 'Avg ...' summary sheets
# ============================================================
# --- Avg Points ---
$ws = $wb.Worksheets.Item("Avg Points")
$ws.Range("A2").Value = "LaMelo Ball"
$ws.Range("B2").Value = 23.33333333333333
$ws.Range("A3").Value = "Miles Bridges"
$ws.Range("B3").Value = 21.81818181818182
$ws.Range("A4").Value = "Kon Knueppel"
$ws.Range("B4").Value = 16.63636363636364
$ws.Range("A5").Value = "Collin Sexton"
$ws.Range("B5").Value = 15.8
$ws.Range("A6").Value = "Brandon Miller"
$ws.Range("B6").Value = 14.5
$ws.Range("A7").Value = "Tre Mann"
$ws.Range("B7").Value = 10.81818181818182
$ws.Range("A8").Value = "Ryan Kalkbrenner"
$ws.Range("B8").Value = 9.909090909090908
$ws.Range("A9").Value = "Moussa Diabaté"
$ws.Range("B9").Value = 9.818181818181818
$ws.Range("A10").Value = "Sion James"
$ws.Range("B10").Value = 8.272727272727273
$ws.Range("A11").Value = "KJ Simpson"
$ws.Range("B11").Value = 5
$ws.Range("A12").Value = "Pat Connaughton"
$ws.Range("B12").Value = 4.285714285714286
$ws.Range("A13").Value = "Liam McNeeley"
$ws.Range("B13").Value = 4.111111111111111
$ws.Range("A14").Value = "Tidjane Salaün"
$ws.Range("B14").Value = 3.142857142857143
$ws.Range("A15").Value = "Mason Plumlee"
$ws.Range("B15").Value = 0.5

# --- Avg Assists ---
$ws = $wb.Worksheets.Item("Avg Assists")
$ws.Range("A2").Value = "LaMelo Ball"
$ws.Range("B2").Value = 9.833333333333334
$ws.Range("A3").Value = "Collin Sexton"
$ws.Range("B3").Value = 5.3
$ws.Range("A4").Value = "Miles Bridges"
$ws.Range("B4").Value = 4.181818181818182
$ws.Range("A5").Value = "Brandon Miller"
$ws.Range("B5").Value = 4
$ws.Range("A6").Value = "Tre Mann"
$ws.Range("B6").Value = 3.181818181818182
$ws.Range("A7").Value = "Kon Knueppel"
$ws.Range("B7").Value = 2.818181818181818
$ws.Range("A8").Value = "KJ Simpson"
$ws.Range("B8").Value = 2.6
$ws.Range("A9").Value = "Sion James"
$ws.Range("B9").Value = 1.363636363636364
$ws.Range("A10").Value = "Liam McNeeley"
$ws.Range("B10").Value = 1.222222222222222
$ws.Range("A11").Value = "Moussa Diabaté"
$ws.Range("B11").Value = 0.7272727272727273
$ws.Range("A12").Value = "Tidjane Salaün"
$ws.Range("B12").Value = 0.7142857142857143
$ws.Range("A13").Value = "Pat Connaughton"
$ws.Range("B13").Value = 0.5714285714285714
$ws.Range("A14").Value = "Ryan Kalkbrenner"
$ws.Range("B14").Value = 0.4545454545454545
$ws.Range("A15").Value = "Mason Plumlee"
$ws.Range("B15").Value = 0.25

# --- Avg Rebounds ---
$ws = $wb.Worksheets.Item("Avg Rebounds")
$ws.Range("A2").Value = "LaMelo Ball"
$ws.Range("B2").Value = 7.833333333333333
$ws.Range("A3").Value = "Moussa Diabaté"
$ws.Range("B3").Value = 7.454545454545454
$ws.Range("A4").Value = "Miles Bridges"
$ws.Range("B4").Value = 7.181818181818182
$ws.Range("A5").Value = "Ryan Kalkbrenner"
$ws.Range("B5").Value = 6.818181818181818
$ws.Range("A6").Value = "Kon Knueppel"
$ws.Range("B6").Value = 6.363636363636363
$ws.Range("A7").Value = "Tidjane Salaün"
$ws.Range("B7").Value = 3.571428571428572
$ws.Range("A8").Value = "Tre Mann"
$ws.Range("B8").Value = 3.090909090909091
$ws.Range("A9").Value = "Sion James"
$ws.Range("B9").Value = 2.727272727272727
$ws.Range("A10").Value = "KJ Simpson"
$ws.Range("B10").Value = 2.6
$ws.Range("A11").Value = "Collin Sexton"
$ws.Range("B11").Value = 2.3
$ws.Range("A12").Value = "Pat Connaughton"
$ws.Range("B12").Value = 2
$ws.Range("A13").Value = "Liam McNeeley"
$ws.Range("B13").Value = 1.888888888888889
$ws.Range("A14").Value = "Mason Plumlee"
$ws.Range("B14").Value = 1
$ws.Range("A15").Value = "Brandon Miller"
$ws.Range("B15").Value = 0

# --- Avg 3PM ---
$ws = $wb.Worksheets.Item("Avg 3PM")
$ws.Range("A2").Value = "LaMelo Ball"
$ws.Range("B2").Value = 3.333333333333333
$ws.Range("A3").Value = "Kon Knueppel"
$ws.Range("B3").Value = 3.181818181818182
$ws.Range("A4").Value = "Miles Bridges"
$ws.Range("B4").Value = 3.090909090909091
$ws.Range("A5").Value = "Tre Mann"
$ws.Range("B5").Value = 1.727272727272727
$ws.Range("A6").Value = "Sion James"
$ws.Range("B6").Value = 1.545454545454545
$ws.Range("A7").Value = "Collin Sexton"
$ws.Range("B7").Value = 1.4
$ws.Range("A8").Value = "Brandon Miller"
$ws.Range("B8").Value = 1
$ws.Range("A9").Value = "KJ Simpson"
$ws.Range("B9").Value = 0.8
$ws.Range("A10").Value = "Pat Connaughton"
$ws.Range("B10").Value = 0.7142857142857143
$ws.Range("A11").Value = "Liam McNeeley"
$ws.Range("B11").Value = 0.6666666666666666
$ws.Range("A12").Value = "Tidjane Salaün"
$ws.Range("B12").Value = 0.5714285714285714
$ws.Range("A13").Value = "Ryan Kalkbrenner"
$ws.Range("B13").Value = 0
$ws.Range("A14").Value = "Moussa Diabaté"
$ws.Range("B14").Value = 0
$ws.Range("A15").Value = "Mason Plumlee"
$ws.Range("B15").Value = 0

